$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by index) whose values move as part of the row rotation.
$cols = @(1,2,4,5,6,7,8,9,10,16,17,18)  # A,B,D,E,F,G,H,I,J,P,Q,R

function Get-RowData($row) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Cells.Item($row, $col).Value2
    }
    return $data
}

# Snapshot the current contents of the affected rows before writing anything,
# since the rotation reassigns rows in terms of each other.
$r8  = Get-RowData 8
$r9  = Get-RowData 9
$r10 = Get-RowData 10
$r11 = Get-RowData 11
$r12 = Get-RowData 12

function Set-RowData($row, $oldData, $newData) {
    foreach ($col in $cols) {
        $newVal = $newData[$col]
        $oldVal = $oldData[$col]
        if ($newVal -eq $null) { $newVal = "" }
        if ($oldVal -eq $null) { $oldVal = "" }
        # Skip writing cells whose value does not actually change, so cells
        # that are empty both before and after are left completely untouched.
        if ($newVal -ne $oldVal) {
            $cell = $ws.Cells.Item($row, $col)
            # The "Antal" column (I) can hold a numeric-looking string (e.g.
            # "1") that must stay text, not become a number - force text
            # format before assigning so Excel doesn't auto-convert it.
            if ($col -eq 9 -and $newVal -ne "") {
                $cell.NumberFormat = "@"
            }
            $cell.Value2 = $newVal
        }
    }
}

# Rotation: row8 <- row10, row9 <- row8, row10 <- row9 ; row11 <-> row12
Set-RowData 8  $r8  $r10
Set-RowData 9  $r9  $r8
Set-RowData 10 $r10 $r9
Set-RowData 11 $r11 $r12
Set-RowData 12 $r12 $r11
